# Update the XMLFieldPath (column C) entries for the "notifications" block so that
# the region-specific "d1:notificationXX" XPath segments are replaced with the
# generic "/*/*" wildcard segments used when running through another region's XML.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C19").Value = "/*/*/oos:id"
$ws.Range("C20").Value = "/*/*/oos:notificationNumber"
$ws.Range("C21").Value = "/*/*/oos:versionNumber"
$ws.Range("C22").Value = "/*/*/oos:publishDate"
$ws.Range("C23").Value = "/*/*/oos:placingWay/oos:name"
$ws.Range("C24").Value = "/*/*/oos:orderName"
$ws.Range("C25").Value = "/*/*/oos:order/oos:placer/oos:regNum"
$ws.Range("C26").Value = "/*/*/oos:order/oos:placer/oos:fullName"

# The old, region-specific notificationZK path is now recorded as an example
# path on the (previously empty) row 28.
$ws.Range("C28").Value = "//*/d1:notificationZK/oos:id"

# Move the active selection to C27, matching the author's last-touched cell.
$ws.Range("C27").Select() | Out-Null
